# Auto-assembled data rows for the ランサーズ sheet (案件情報.xlsx)
# Rebuilds rows 2-15 to match the new scraped/sorted dataset, updates the
# dimension, and recreates hyperlinks for column F with the Hyperlink style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all existing hyperlinks first -- Hyperlinks.Add() stacks a new link on
# top of an existing one rather than replacing it, and every F-cell in this
# sheet is being rewritten with a (possibly new) URL anyway.
$ws.Hyperlinks.Delete()

# --- Row data: Timestamp, Title, Category, Price, Deadline, URL, Score, Skill ---
$rows = @()
$rows += ,@(2, "2026-01-19 18:29:12", "急募】Google API + LINE API + OPEN AI API /Web SaaS開発", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474570", 535, "🔥AI,API ◆開発")
$rows += ,@(3, "2026-01-19 18:29:12", "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5423720", 385, "🔥AI,Ai ◆効率化")
$rows += ,@(4, "2026-01-19 18:29:12", "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5434128", 368, "🔥AI,Ai ◆開発")
$rows += ,@(5, "2026-01-19 18:29:12", "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5427956", 310, "🔥AI,Ai")
$rows += ,@(6, "2026-01-19 18:29:12", "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5439158", 303, "🔥AI,Ai")
$rows += ,@(7, "2026-01-19 18:29:12", "【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,300円前後)", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474415", 213, "🔥API ◇管理")
$rows += ,@(8, "2026-01-19 18:29:12", "【急募】外国人社員のシフト・欠勤管理アプリ開発依頼", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474622", 135, "◆開発 ◇アプリ")
$rows += ,@(9, "2026-01-19 18:29:12", "【医療機関向け業務改善サービスの新規開発】WEBアプリ開発におけるフルスタック開発担当者募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5473940", 135, "◆開発 ◇業務改善")
$rows += ,@(10, "2026-01-19 18:29:12", "【急募】インバウンド向け新サービスアプリ開発見積作成依頼", "システム開発", "1,000 ~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474189", 85, "◆開発 ◇アプリ")
$rows += ,@(11, "2026-01-19 18:29:12", "進行管理およびチームディレクションを担当", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5418064", 30, "◇管理")
$rows += ,@(12, "2026-01-19 18:29:12", "JSを使用したSaaSサービスの導入、保守のパートナー募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474214", 25, $null)
$rows += ,@(13, "2026-01-19 18:29:12", "【急募】自動車整備業向けCRM構築パートナー募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474125", 25, $null)
$rows += ,@(14, "2026-01-19 18:29:12", "移動型演出カートの電装設計および制御ユニット製作(Arduino/ESP32等)", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474273", 18, $null)
$rows += ,@(15, "2026-01-19 18:29:12", "Google clab用マークシートCSV出力プログラム作成依頼", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5474679", 10, $null)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]

    $urlCell = $ws.Cells.Item($r, 6)
    $urlCell.Value = $row[6]
    $ws.Hyperlinks.Add($urlCell, $row[6])
    $urlCell.Style = "Hyperlink"

    $ws.Cells.Item($r, 7).Value = $row[7]

    $skill = $row[8]
    if ($skill -eq $null) {
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $skill
    }
}

Write-Output "done"
